# SIM800L SMS & Paypal update
# - adds smarthome info (row 8)
# - renames/relabels coin power -> coin power relais (row 19)
# - adds Buzzer / Paypal info (row 51)
# - adds Sim800L SMS module TX/RX rows (53, 55)
# - removes old "Open Box" note (row 54)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mega 2560")

# Row 8: new smarthome relay note
$ws.Range("C8").Value = "smarthome"
$ws.Range("E8").Value = "smarthome info about sold items"

# Row 19: "coin power" -> "coin power`n relais", now wrapped (style already used by E19)
$ws.Range("C19").Value = "coin power" + [char]10 + " relais"
$ws.Range("C19").WrapText = $true

# Row 51: Buzzer to inform about incoming Paypal money transfer (clear old style first)
$ws.Range("C51").ClearFormats()
$ws.Range("E51").ClearFormats()
$ws.Range("C51").Value = "Buzzer"
$ws.Range("E51").Value = "Buzzer to inform about incoming Paypal money transfer"

# Row 53: Sim800L SMS module - TX pin
$ws.Range("C53").Value = "TX"
$ws.Range("E53").Value = "Sim 800L SMS module "

# Row 54: remove old "Open Box - if electronical lock is implemented" / white note
$ws.Range("E54:F54").ClearContents()

# Row 55: Sim800L SMS module - RX pin (replaces old smarthome note text)
$ws.Range("C55").Value = "RX"
$ws.Range("E55").Value = "Sim 800L SMS module "

# Update selection / navigation state to mirror the saved view
$ws.Range("I60").Select()
